$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.940.03"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.91"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.16"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5077"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06385"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.303"
$ws.Range("E12").Value = "  +0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.93"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5467"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7863"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.52"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.985.88"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.03"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.432"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.963"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.051"
$ws.Range("E22").Value = "  +1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.877"
$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.44"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1143"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.882"
$ws.Range("E27").Value = "  +2.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.74"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05045"
$ws.Range("E29").Value = "  +2.10%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.237"
$ws.Range("E30").Value = "  -0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.260"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.541"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.362"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8948"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.597"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134.45"
$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5492"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("B39").Value = "BabyDogeCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₈133"
$ws.Range("E39").Value = "  +14.84%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.540"
$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.626"
$ws.Range("E43").Value = "  -1.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8157"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.779.45"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.90"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05072"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  +0.36%  "
